$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 1: ALC!row46 (G46=4584)
$ws.Range("H46").Value = 2433.3333
$ws.Range("I46").Value = 1150
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 3450
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -3331
$ws.Range("N46").Value = -15238

# Hunk 2: ALC!row60 (G60=4584)
$ws.Range("H60").Value = 2433.3333
$ws.Range("I60").Value = 1150
$ws.Range("J60").Value = 5000
$ws.Range("K60").Value = 3450
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = -2966
$ws.Range("N60").Value = -15968

# Hunk 3: ALC!row125 (G125=36228)
$ws.Range("H125").Value = 1406.2
$ws.Range("I125").Value = 1066
$ws.Range("K125").Value = 9594
$ws.Range("M125").Value = -7134

$ws = $wb.Worksheets.Item("ARM")
# Hunk 4: ARM!row2 (G2=27713)
$ws.Range("H2").Value = 1325.8667
$ws.Range("I2").Value = 1299.0769
$ws.Range("K2").Value = 1299.0769
$ws.Range("M2").Value = -1186.0769

# Hunk 5: ARM!row32 (G32=44147)
$ws.Range("H32").Value = 6469.4814
$ws.Range("I32").Value = 5711.7026
$ws.Range("J32").Value = 8118.7646
$ws.Range("K32").Value = 5711.7026
$ws.Range("L32").Value = 8118.7646
$ws.Range("M32").Value = -5424.7026
$ws.Range("N32").Value = -8692.7646

# Hunk 6: ARM!row74 (G74=44000)
$ws.Range("H74").Value = 7910.154
$ws.Range("I74").Value = 8381.799999999999
$ws.Range("K74").Value = 8381.799999999999
$ws.Range("M74").Value = -7507.799999999999

# Hunk 7: ARM!row76 (G76=10679)
$ws.Range("H76").Value = 37679.8
$ws.Range("J76").Value = 37679.8
$ws.Range("L76").Value = 37679.8
$ws.Range("N76").Value = -38355.8

# Hunk 8: ARM!row77 (G77=44000)
$ws.Range("H77").Value = 7910.154
$ws.Range("I77").Value = 8381.799999999999
$ws.Range("K77").Value = 41909
$ws.Range("M77").Value = -37541

# Hunk 9: ARM!row79 (G79=10679)
$ws.Range("H79").Value = 37679.8
$ws.Range("J79").Value = 37679.8
$ws.Range("L79").Value = 37679.8
$ws.Range("N79").Value = -40019.8

# Hunk 10: ARM!row110 (G110=27708)
$ws.Range("H110").Value = 863.6061
$ws.Range("I110").Value = 864.70966
$ws.Range("J110").Value = 846.5
$ws.Range("K110").Value = 864.70966
$ws.Range("L110").Value = 846.5
$ws.Range("M110").Value = 1180.29034
$ws.Range("N110").Value = -4936.5

# Hunk 11: ARM!row116 (G116=27713)
$ws.Range("H116").Value = 1325.8667
$ws.Range("I116").Value = 1299.0769
$ws.Range("K116").Value = 1299.0769
$ws.Range("M116").Value = 994.9231

$ws = $wb.Worksheets.Item("BSM")
# Hunk 12: BSM!row3 (G3=27713)
$ws.Range("H3").Value = 1325.8667
$ws.Range("I3").Value = 1299.0769
$ws.Range("K3").Value = 1299.0769
$ws.Range("M3").Value = -1185.0769

# Hunk 13: BSM!row105 (G105=19947)
$ws.Range("H105").Value = 1601.4667
$ws.Range("J105").Value = 1397.3334
$ws.Range("L105").Value = 1397.3334
$ws.Range("N105").Value = -4891.3334

# Hunk 14: BSM!row118 (G118=27137)
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# Hunk 15: BSM!row134 (G134=43998)
$ws.Range("H134").Value = 1999.2927
$ws.Range("I134").Value = 1539.7567
$ws.Range("J134").Value = 6250
$ws.Range("K134").Value = 4619.2701
$ws.Range("L134").Value = 18750
$ws.Range("M134").Value = -2084.2701
$ws.Range("N134").Value = -23820

$ws = $wb.Worksheets.Item("CRP")
# Hunk 16: CRP!row16 (G16=27691)
$ws.Range("H16").Value = 6945579
$ws.Range("I16").Value = 9260272
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 9260272
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -9259985
$ws.Range("N16").Value = -2074

# Hunk 17: CRP!row31 (G31=44023)
$ws.Range("H31").Value = 17244476
$ws.Range("I31").Value = 1466.2667
$ws.Range("J31").Value = 35719130
$ws.Range("K31").Value = 1466.2667
$ws.Range("L31").Value = 35719130
$ws.Range("M31").Value = -1171.2667
$ws.Range("N31").Value = -35719720

# Hunk 18: CRP!row34 (G34=44023)
$ws.Range("H34").Value = 17244476
$ws.Range("I34").Value = 1466.2667
$ws.Range("J34").Value = 35719130
$ws.Range("K34").Value = 1466.2667
$ws.Range("L34").Value = 35719130
$ws.Range("M34").Value = -1264.2667
$ws.Range("N34").Value = -35719534

# Hunk 19: CRP!row58 (G58=44021)
$ws.Range("H58").Value = 1543.8914
$ws.Range("I58").Value = 1332.4097
$ws.Range("J58").Value = 3494.2222
$ws.Range("K58").Value = 1332.4097
$ws.Range("L58").Value = 3494.2222
$ws.Range("M58").Value = -1129.4097
$ws.Range("N58").Value = -3900.2222

# Hunk 20: CRP!row81 (G81=10613)
$ws.Range("H81").Value = 27966.334
$ws.Range("J81").Value = 27966.334
$ws.Range("L81").Value = 27966.334
$ws.Range("N81").Value = -29962.334

# Hunk 21: CRP!row84 (G84=10613)
$ws.Range("H84").Value = 27966.334
$ws.Range("J84").Value = 27966.334
$ws.Range("L84").Value = 83899.00199999999
$ws.Range("N84").Value = -93883.00199999999

# Hunk 22: CRP!row113 (G113=27691)
$ws.Range("H113").Value = 6945579
$ws.Range("I113").Value = 9260272
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 9260272
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -9258102
$ws.Range("N113").Value = -5840

# Hunk 23: CRP!row136 (G136=44021)
$ws.Range("H136").Value = 1543.8914
$ws.Range("I136").Value = 1332.4097
$ws.Range("J136").Value = 3494.2222
$ws.Range("K136").Value = 3997.2291
$ws.Range("L136").Value = 10482.6666
$ws.Range("M136").Value = -1447.2291
$ws.Range("N136").Value = -15582.6666

$ws = $wb.Worksheets.Item("CUL")
# Hunk 24: CUL!row80 (G80=12890)
$ws.Range("H80").Value = 9998.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 9998.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 29996.25
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -31868.25

# Hunk 25: CUL!row83 (G83=12890)
$ws.Range("H83").Value = 9998.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 9998.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 89988.75
$ws.Range("M83").ClearContents()  # becomes empty-valued cell (approximation)
$ws.Range("N83").Value = -99348.75

# Hunk 26: CUL!row131 (G131=36060)
$ws.Range("H131").Value = 7693196.5
$ws.Range("J131").Value = 919.2787
$ws.Range("L131").Value = 2757.8361
$ws.Range("N131").Value = -12837.8361

$ws = $wb.Worksheets.Item("GSM")
# Hunk 27: GSM!row121 (G121=26338)
$ws.Range("H121").Value = 27645
$ws.Range("J121").Value = 27645
$ws.Range("L121").Value = 27645
$ws.Range("N121").Value = -31139

$ws = $wb.Worksheets.Item("LTW")
# Hunk 28: LTW!row74 (G74=11990)
$ws.Range("H74").Value = 45400
$ws.Range("J74").Value = 45400
$ws.Range("L74").Value = 45400
$ws.Range("N74").Value = -47396

# Hunk 29: LTW!row77 (G77=11990)
$ws.Range("H77").Value = 45400
$ws.Range("J77").Value = 45400
$ws.Range("L77").Value = 136200
$ws.Range("N77").Value = -146184

$ws = $wb.Worksheets.Item("WVR")
# Hunk 30: WVR!row132 (G132=44029)
$ws.Range("H132").Value = 10102736
$ws.Range("I132").Value = 1054.9615
$ws.Range("J132").Value = 47623264
$ws.Range("K132").Value = 3164.8845
$ws.Range("L132").Value = 142869792
$ws.Range("M132").Value = -634.8844999999997
$ws.Range("N132").Value = -142874852
Write-Output "Applied 30 cell-group updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
